$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 773.619
$ws.Range("I41").Value = 1297.2
$ws.Range("J41").Value = 297.63635
$ws.Range("K41").Value = 1297.2
$ws.Range("L41").Value = 297.63635
$ws.Range("M41").Value = -857.2
$ws.Range("N41").Value = -1177.63635
$ws.Range("H76").Value = 374887
$ws.Range("I76").Value = 1113611
$ws.Range("J76").Value = 5525
$ws.Range("K76").Value = 1113611
$ws.Range("L76").Value = 5525
$ws.Range("M76").Value = -1113296
$ws.Range("N76").Value = -6155
$ws.Range("H79").Value = 374887
$ws.Range("I79").Value = 1113611
$ws.Range("J79").Value = 5525
$ws.Range("K79").Value = 1113611
$ws.Range("L79").Value = 5525
$ws.Range("M79").Value = -1112519
$ws.Range("N79").Value = -7709
$ws.Range("H86").Value = 10563090
$ws.Range("I86").Value = 12071889
$ws.Range("J86").Value = 1502
$ws.Range("K86").Value = 12071889
$ws.Range("L86").Value = 1502
$ws.Range("M86").Value = -12070766
$ws.Range("N86").Value = -3748
$ws.Range("H89").Value = 10563090
$ws.Range("I89").Value = 12071889
$ws.Range("J89").Value = 1502
$ws.Range("K89").Value = 60359445
$ws.Range("L89").Value = 7510
$ws.Range("M89").Value = -60353829
$ws.Range("N89").Value = -18742
$ws.Range("H96").Value = 1095.1428
$ws.Range("I96").Value = 442.63635
$ws.Range("J96").Value = 1812.9
$ws.Range("K96").Value = 1327.90905
$ws.Range("L96").Value = 5438.700000000001
$ws.Range("M96").Value = 45.09095000000002
$ws.Range("N96").Value = -8184.700000000001
$ws.Range("H98").Value = 750
$ws.Range("I98").Value = 733.3333
$ws.Range("J98").Value = 800
$ws.Range("K98").Value = 733.3333
$ws.Range("L98").Value = 800
$ws.Range("M98").Value = 764.6667
$ws.Range("N98").Value = -3796
$ws.Range("H122").Value = 750
$ws.Range("I122").Value = 733.3333
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 2199.9999
$ws.Range("L122").Value = 2400
$ws.Range("M122").Value = 250.0001000000002
$ws.Range("N122").Value = -7300
$ws.Range("H129").Value = 1026.8
$ws.Range("J129").Value = 1313.6364
$ws.Range("L129").Value = 3940.9092
$ws.Range("N129").Value = -13940.9092
$ws.Range("H132").Value = 2959.239
$ws.Range("I132").Value = 2958.3333
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 8874.999899999999
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -6344.999899999999
$ws.Range("N132").Value = -14060
$ws.Range("H137").Value = 1820.7084
$ws.Range("I137").Value = 1543.6428
$ws.Range("J137").Value = 2208.6
$ws.Range("K137").Value = 4630.928400000001
$ws.Range("L137").Value = 6625.799999999999
$ws.Range("M137").Value = -2080.928400000001
$ws.Range("N137").Value = -11725.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1390.85
$ws.Range("I45").Value = 1238.4546
$ws.Range("J45").Value = 1577.1111
$ws.Range("K45").Value = 1238.4546
$ws.Range("L45").Value = 1577.1111
$ws.Range("M45").Value = -861.4546
$ws.Range("N45").Value = -2331.1111
$ws.Range("H61").Value = 1925.7693
$ws.Range("I61").Value = 1661.8334
$ws.Range("K61").Value = 1661.8334
$ws.Range("M61").Value = -1449.8334
$ws.Range("H132").Value = 1485.0526
$ws.Range("I132").Value = 1257.5686
$ws.Range("J132").Value = 3418.6667
$ws.Range("K132").Value = 3772.7058
$ws.Range("L132").Value = 10256.0001
$ws.Range("M132").Value = -1242.7058
$ws.Range("N132").Value = -15316.0001
$ws.Range("H136").Value = 1925.7693
$ws.Range("I136").Value = 1661.8334
$ws.Range("K136").Value = 4985.5002
$ws.Range("M136").Value = -2435.5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1874.4736
$ws.Range("I105").Value = 1677.0834
$ws.Range("J105").Value = 2212.8572
$ws.Range("K105").Value = 1677.0834
$ws.Range("L105").Value = 2212.8572
$ws.Range("M105").Value = 69.91660000000002
$ws.Range("N105").Value = -5706.8572
$ws.Range("H107").Value = 1728.8948
$ws.Range("I107").Value = 1527.7858
$ws.Range("J107").Value = 2292
$ws.Range("K107").Value = 1527.7858
$ws.Range("L107").Value = 2292
$ws.Range("M107").Value = 392.2141999999999
$ws.Range("N107").Value = -6132
$ws.Range("H134").Value = 2686.625
$ws.Range("I134").Value = 1675.8
$ws.Range("J134").Value = 4371.3335
$ws.Range("K134").Value = 5027.4
$ws.Range("L134").Value = 13114.0005
$ws.Range("M134").Value = -2492.4
$ws.Range("N134").Value = -18184.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1770.6957
$ws.Range("I16").Value = 1769.375
$ws.Range("J16").Value = 1773.7142
$ws.Range("K16").Value = 1769.375
$ws.Range("L16").Value = 1773.7142
$ws.Range("M16").Value = -1482.375
$ws.Range("N16").Value = -2347.7142
$ws.Range("H31").Value = 13368311
$ws.Range("I31").Value = 10580070
$ws.Range("K31").Value = 10580070
$ws.Range("M31").Value = -10579775
$ws.Range("H34").Value = 13368311
$ws.Range("I34").Value = 10580070
$ws.Range("K34").Value = 10580070
$ws.Range("M34").Value = -10579868
$ws.Range("H113").Value = 1770.6957
$ws.Range("I113").Value = 1769.375
$ws.Range("J113").Value = 1773.7142
$ws.Range("K113").Value = 1769.375
$ws.Range("L113").Value = 1773.7142
$ws.Range("M113").Value = 400.625
$ws.Range("N113").Value = -6113.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 17.444445
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 18.411764
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 55.235292
$ws.Range("M12").Value = 170
$ws.Range("N12").Value = -401.235292
$ws.Range("H45").Value = 806
$ws.Range("I45").Value = 676.6667
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 2030.0001
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = -1498.0001
$ws.Range("N45").Value = -4064
$ws.Range("H127").Value = 1933
$ws.Range("J127").Value = 1933
$ws.Range("L127").Value = 5799
$ws.Range("N127").Value = -15719
$ws.Range("H130").Value = 1332.5
$ws.Range("I130").Value = 1115
$ws.Range("J130").Value = 1550
$ws.Range("K130").Value = 3345
$ws.Range("L130").Value = 4650
$ws.Range("M130").Value = 1675
$ws.Range("N130").Value = -14690
$ws.Range("H131").Value = 1361385.9
$ws.Range("I131").Value = 6060937.5
$ws.Range("J131").Value = 989.3421
$ws.Range("K131").Value = 18182812.5
$ws.Range("L131").Value = 2968.0263
$ws.Range("M131").Value = -18177772.5
$ws.Range("N131").Value = -13048.0263
$ws.Range("H133").Value = 1896.1111
$ws.Range("I133").Value = 1552.5
$ws.Range("J133").Value = 1994.2858
$ws.Range("K133").Value = 4657.5
$ws.Range("L133").Value = 5982.857400000001
$ws.Range("M133").Value = 402.5
$ws.Range("N133").Value = -16102.8574
$ws.Range("H139").Value = 1462.826
$ws.Range("I139").Value = 937.4545000000001
$ws.Range("J139").Value = 1944.4166
$ws.Range("K139").Value = 2812.3635
$ws.Range("L139").Value = 5833.2498
$ws.Range("M139").Value = 2327.6365
$ws.Range("N139").Value = -16113.2498

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1915.8572
$ws.Range("I113").Value = 1637
$ws.Range("J113").Value = 2125
$ws.Range("K113").Value = 1637
$ws.Range("L113").Value = 2125
$ws.Range("M113").Value = 533
$ws.Range("N113").Value = -6465
$ws.Range("H132").Value = 2094.5862
$ws.Range("I132").Value = 1596.9565
$ws.Range("J132").Value = 4002.1667
$ws.Range("K132").Value = 4790.8695
$ws.Range("L132").Value = 12006.5001
$ws.Range("M132").Value = -2260.8695
$ws.Range("N132").Value = -17066.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 917.3333
$ws.Range("I16").Value = 907
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 907
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -737
$ws.Range("N16").Value = -1340
$ws.Range("H61").Value = 1861.25
$ws.Range("I61").Value = 1445
$ws.Range("J61").Value = 2277.5
$ws.Range("K61").Value = 1445
$ws.Range("L61").Value = 2277.5
$ws.Range("M61").Value = -1243
$ws.Range("N61").Value = -2681.5
$ws.Range("H113").Value = 1861.25
$ws.Range("I113").Value = 1445
$ws.Range("J113").Value = 2277.5
$ws.Range("K113").Value = 1445
$ws.Range("L113").Value = 2277.5
$ws.Range("M113").Value = 725
$ws.Range("N113").Value = -6617.5
$ws.Range("H132").Value = 2853418.8
$ws.Range("I132").Value = 4481430
$ws.Range("J132").Value = 4399.375
$ws.Range("K132").Value = 13444290
$ws.Range("L132").Value = 13198.125
$ws.Range("M132").Value = -13441760
$ws.Range("N132").Value = -18258.125
$ws.Range("H136").Value = 5006759.5
$ws.Range("I136").Value = 6586826
$ws.Range("J136").Value = 3216.6667
$ws.Range("K136").Value = 19760478
$ws.Range("L136").Value = 9650.000100000001
$ws.Range("M136").Value = -19757928
$ws.Range("N136").Value = -14750.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 62501370
$ws.Range("I113").Value = 71430380
$ws.Range("J113").Value = 55556576
$ws.Range("K113").Value = 214291140
$ws.Range("L113").Value = 166669728
$ws.Range("M113").Value = -214288970
$ws.Range("N113").Value = -166674068
$ws.Range("H132").Value = 2450.4707
$ws.Range("I132").Value = 1968.8572
$ws.Range("J132").Value = 4698
$ws.Range("K132").Value = 5906.571599999999
$ws.Range("L132").Value = 14094
$ws.Range("M132").Value = -3376.571599999999
$ws.Range("N132").Value = -19154
$ws.Range("H136").Value = 1676.875
$ws.Range("I136").Value = 1474.75
$ws.Range("J136").Value = 1879
$ws.Range("K136").Value = 4424.25
$ws.Range("L136").Value = 5637
$ws.Range("M136").Value = -1874.25
$ws.Range("N136").Value = -10737

Write-Host "Applied profit/price updates to 255 cells across 8 sheets"
